$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 121
$ws.Range("H121").Value = 1397.1666
$ws.Range("I121").Value = 641.5
$ws.Range("J121").Value = 1775
$ws.Range("K121").Value = 1924.5
$ws.Range("L121").Value = 5325
$ws.Range("M121").Value = -177.5
$ws.Range("N121").Value = -8819
# Row 137
$ws.Range("H137").Value = 9098455
$ws.Range("I137").Value = 25015500
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 75046500
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -75043950
$ws.Range("N137").Value = -14100
# Row 138
$ws.Range("H138").Value = 3953.606
$ws.Range("I138").Value = 2257.0588
$ws.Range("J138").Value = 5756.1875
$ws.Range("K138").Value = 6771.176399999999
$ws.Range("L138").Value = 17268.5625
$ws.Range("M138").Value = -1631.176399999999
$ws.Range("N138").Value = -27548.5625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1861.5385
$ws.Range("J45").Value = 6075
$ws.Range("L45").Value = 6075
$ws.Range("N45").Value = -6829
# Row 58
$ws.Range("H58").Value = 28500
$ws.Range("J58").Value = 28500
$ws.Range("L58").Value = 28500
$ws.Range("N58").Value = -29360
# Row 97
$ws.Range("H97").Value = 521.23334
$ws.Range("I97").Value = 513.2857
$ws.Range("J97").Value = 632.5
$ws.Range("K97").Value = 513.2857
$ws.Range("L97").Value = 632.5
$ws.Range("M97").Value = -17.28570000000002
$ws.Range("N97").Value = -1624.5
# Row 122
$ws.Range("H122").Value = 3453.6667
$ws.Range("I122").Value = 2240.6667
$ws.Range("J122").Value = 4666.6665
$ws.Range("K122").Value = 6722.000100000001
$ws.Range("L122").Value = 13999.9995
$ws.Range("M122").Value = -4272.000100000001
$ws.Range("N122").Value = -18899.9995
# Row 132
$ws.Range("H132").Value = 47624850
$ws.Range("I132").Value = 71434984
$ws.Range("J132").Value = 4575.4287
$ws.Range("K132").Value = 214304952
$ws.Range("L132").Value = 13726.2861
$ws.Range("M132").Value = -214302422
$ws.Range("N132").Value = -18786.2861

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3033
$ws.Range("I134").Value = 1918.3125
$ws.Range("J134").Value = 6600
$ws.Range("K134").Value = 5754.9375
$ws.Range("L134").Value = 19800
$ws.Range("M134").Value = -3219.9375
$ws.Range("N134").Value = -24870
# Row 135
$ws.Range("H135").Value = 31580
$ws.Range("J135").Value = 31580
$ws.Range("L135").Value = 31580
$ws.Range("N135").Value = -41720

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1925886.2
$ws.Range("I31").Value = 2327490.2
$ws.Range("K31").Value = 2327490.2
$ws.Range("M31").Value = -2327195.2
# Row 34
$ws.Range("H34").Value = 1925886.2
$ws.Range("I34").Value = 2327490.2
$ws.Range("K34").Value = 2327490.2
$ws.Range("M34").Value = -2327288.2
# Row 41
$ws.Range("H41").Value = 4961.625
$ws.Range("J41").Value = 8450
$ws.Range("L41").Value = 8450
$ws.Range("N41").Value = -9306
# Row 50
$ws.Range("H50").Value = 11066.667
$ws.Range("I50").Value = 1200
$ws.Range("J50").Value = 16000
$ws.Range("K50").Value = 1200
$ws.Range("L50").Value = 16000
$ws.Range("M50").Value = -575
$ws.Range("N50").Value = -17250
# Row 51
$ws.Range("H51").Value = 200013180
$ws.Range("J51").Value = 16475
$ws.Range("L51").Value = 16475
$ws.Range("N51").Value = -17947
# Row 59
$ws.Range("H59").Value = 18500
$ws.Range("J59").Value = 18500
$ws.Range("L59").Value = 18500
$ws.Range("N59").Value = -20790
# Row 60
$ws.Range("H60").Value = 11867.777
$ws.Range("J60").Value = 11697.2
$ws.Range("L60").Value = 11697.2
$ws.Range("N60").Value = -12719.2
# Row 61
$ws.Range("H61").Value = 200013180
$ws.Range("J61").Value = 16475
$ws.Range("L61").Value = 16475
$ws.Range("N61").Value = -17171
# Row 68
$ws.Range("H68").Value = 46000
$ws.Range("J68").Value = 46000
$ws.Range("L68").Value = 46000
$ws.Range("N68").Value = -47498
# Row 71
$ws.Range("H71").Value = 46000
$ws.Range("J71").Value = 46000
$ws.Range("L71").Value = 138000
$ws.Range("N71").Value = -145488
# Row 74
$ws.Range("H74").Value = 20054.223
$ws.Range("J74").Value = 20054.223
$ws.Range("L74").Value = 20054.223
$ws.Range("N74").Value = -21802.223
# Row 77
$ws.Range("H77").Value = 20054.223
$ws.Range("J77").Value = 20054.223
$ws.Range("L77").Value = 60162.66900000001
$ws.Range("N77").Value = -68898.66900000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 692.38635
$ws.Range("I113").Value = 534
$ws.Range("J113").Value = 865.8570999999999
$ws.Range("K113").Value = 1602
$ws.Range("L113").Value = 2597.5713
$ws.Range("M113").Value = 568
$ws.Range("N113").Value = -6937.5713

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 6200.4707
$ws.Range("I122").Value = 9350
$ws.Range("J122").Value = 5231.385
$ws.Range("K122").Value = 28050
$ws.Range("L122").Value = 15694.155
$ws.Range("M122").Value = -25600
$ws.Range("N122").Value = -20594.155
# Row 132
$ws.Range("H132").Value = 2974.484
$ws.Range("I132").Value = 2567.2778
$ws.Range("J132").Value = 3538.3076
$ws.Range("K132").Value = 7701.8334
$ws.Range("L132").Value = 10614.9228
$ws.Range("M132").Value = -5171.8334
$ws.Range("N132").Value = -15674.9228
# Row 139
$ws.Range("H139").Value = 28137.666
$ws.Range("J139").Value = 28137.666
$ws.Range("L139").Value = 28137.666
$ws.Range("N139").Value = -38417.666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 2656.25
$ws.Range("I122").Value = 2462.963
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 7388.889000000001
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -4938.889000000001
$ws.Range("N122").Value = -16000
# Row 132
$ws.Range("H132").Value = 1924.75
$ws.Range("I132").Value = 1182.3478
$ws.Range("J132").Value = 3822
$ws.Range("K132").Value = 3547.0434
$ws.Range("L132").Value = 11466
$ws.Range("M132").Value = -1017.0434
$ws.Range("N132").Value = -16526
# Row 136
$ws.Range("H136").Value = 4002245.2
$ws.Range("I136").Value = 6251284
$ws.Range("J136").Value = 3954.111
$ws.Range("K136").Value = 18753852
$ws.Range("L136").Value = 11862.333
$ws.Range("M136").Value = -18751302
$ws.Range("N136").Value = -16962.333
# Row 141
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("M141").Value = -50360
$ws.Range("N141").Value = -50360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1227.52
$ws.Range("I113").Value = 599.125
$ws.Range("K113").Value = 1797.375
$ws.Range("M113").Value = 372.625
# Row 132
$ws.Range("H132").Value = 332700.44
$ws.Range("I132").Value = 505395.56
$ws.Range("J132").Value = 44875.25
$ws.Range("K132").Value = 1516186.68
$ws.Range("L132").Value = 134625.75
$ws.Range("M132").Value = -1513656.68
$ws.Range("N132").Value = -139685.75
# Row 135
$ws.Range("H135").Value = 68848.75
$ws.Range("J135").Value = 68848.75
$ws.Range("L135").Value = 68848.75
$ws.Range("N135").Value = -78988.75
# Row 136
$ws.Range("H136").Value = 2464
$ws.Range("I136").Value = 1217.3334
$ws.Range("J136").Value = 3960
$ws.Range("K136").Value = 3652.0002
$ws.Range("L136").Value = 11880
$ws.Range("M136").Value = -1102.0002
$ws.Range("N136").Value = -16980
